# "updated activity till excel form" -- the per-innings batting log on the
# "Rohit Sharma (c)" sheet got refreshed: the runs/balls/fours/sixes rows
# were re-ordered to reflect the latest match activity.
#
# The sheet stores these numeric-looking figures as TEXT (the original file
# has t="str" cells, not numeric cells), so we force the number format to
# "@" (Text) before writing each value as a string -- otherwise Excel's
# COM Range.Value setter would auto-coerce a numeric-looking string into a
# real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:F13").NumberFormat = "@"

$ws.Range("C2").Value = "9"
$ws.Range("D2").Value = "8"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "0"

$ws.Range("C3").Value = "35"
$ws.Range("D3").Value = "36"
$ws.Range("E3").Value = "5"
$ws.Range("F3").Value = "1"

$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "7"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"

$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "1"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0"

$ws.Range("C6").Value = "35"
$ws.Range("D6").Value = "23"
$ws.Range("E6").Value = "2"
$ws.Range("F6").Value = "3"

$ws.Range("C7").Value = "5"
$ws.Range("D7").Value = "12"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "12"
$ws.Range("D8").Value = "10"
$ws.Range("E8").Value = "2"
$ws.Range("F8").Value = "0"

$ws.Range("C9").Value = "68"
$ws.Range("D9").Value = "51"
$ws.Range("E9").Value = "5"
$ws.Range("F9").Value = "4"

$ws.Range("C10").Value = "80"
$ws.Range("D10").Value = "54"
$ws.Range("E10").Value = "3"
$ws.Range("F10").Value = "6"

$ws.Range("C11").Value = "6"
$ws.Range("D11").Value = "5"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "1"

$ws.Range("C12").Value = "8"
$ws.Range("D12").Value = "8"
$ws.Range("E12").Value = "0"
$ws.Range("F12").Value = "1"

$ws.Range("C13").Value = "70"
$ws.Range("D13").Value = "45"
$ws.Range("E13").Value = "8"
$ws.Range("F13").Value = "3"
